$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D cells are treated as literal text so values like
# "1.00", "143.80", "48.50" etc. are not normalized into numbers.
$dCells = @("D2", "D3", "D5", "D6", "D9", "D10", "D12", "D13", "D14", "D16", "D17", "D18", "D19", "D20", "D24", "D25", "D27", "D28", "D30", "D33", "D34", "D36", "D37", "D39", "D41", "D44", "D45", "D46", "D47", "D48", "D50")
foreach ($cellRef in $dCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Column D (Price) updates
$ws.Range("D2").Value = '62.457.49'
$ws.Range("D3").Value = '2.898.14'
$ws.Range("D5").Value = '567.35'
$ws.Range("D6").Value = '143.80'
$ws.Range("D9").Value = '2.893.38'
$ws.Range("D10").Value = '6.88'
$ws.Range("D12").Value = '0.430'
$ws.Range("D13").Value = '0.0000237'
$ws.Range("D14").Value = '32.73'
$ws.Range("D16").Value = '3.380.56'
$ws.Range("D17").Value = '62.424.29'
$ws.Range("D18").Value = '6.58'
$ws.Range("D19").Value = '2.898.50'
$ws.Range("D20").Value = '425.21'
$ws.Range("D24").Value = '78.58'
$ws.Range("D25").Value = '11.87'
$ws.Range("D27").Value = '1.00'
$ws.Range("D28").Value = '2.01'
$ws.Range("D30").Value = '7.09'
$ws.Range("D33").Value = '0.999'
$ws.Range("D34").Value = '25.67'
$ws.Range("D36").Value = '0.944'
$ws.Range("D37").Value = '5.37'
$ws.Range("D39").Value = '48.50'
$ws.Range("D41").Value = '41.23'
$ws.Range("D44").Value = '0.265'
$ws.Range("D45").Value = '2.714.69'
$ws.Range("D46").Value = '133.18'
$ws.Range("D47").Value = '0.0337'
$ws.Range("D48").Value = '355.03'
$ws.Range("D50").Value = '0.000218'

# Column E (Volume 1h) updates
$ws.Range("E2").Value = '  +0.98%  '
$ws.Range("E3").Value = '  -0.42%  '
$ws.Range("E4").Value = '  +0.12%  '
$ws.Range("E5").Value = '  -3.30%  '
$ws.Range("E6").Value = '  -1.87%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("E8").Value = '  -1.39%  '
$ws.Range("E9").Value = '  -0.50%  '
$ws.Range("E10").Value = '  +0.49%  '
$ws.Range("E11").Value = '  -0.49%  '
$ws.Range("E12").Value = '  -1.09%  '
$ws.Range("E13").Value = '  +0.10%  '
$ws.Range("E14").Value = '  -0.24%  '
$ws.Range("E15").Value = '  +0.24%  '
$ws.Range("E16").Value = '  -0.35%  '
$ws.Range("E17").Value = '  +0.90%  '
$ws.Range("E18").Value = '  -0.93%  '
$ws.Range("E19").Value = '  -0.43%  '
$ws.Range("E20").Value = '  -2.64%  '
$ws.Range("E21").Value = '  -2.73%  '
$ws.Range("E22").Value = '  -0.58%  '
$ws.Range("E23").Value = '  -1.56%  '
$ws.Range("E25").Value = '  -0.84%  '
$ws.Range("E26").Value = '  -2.13%  '
$ws.Range("E27").Value = '  -0.05%  '
$ws.Range("E28").Value = '  -3.20%  '
$ws.Range("E29").Value = '  +1.66%  '
$ws.Range("E30").Value = '  -0.92%  '
$ws.Range("E31").Value = '  -3.15%  '
$ws.Range("E32").Value = '  -4.50%  '
$ws.Range("E33").Value = '  -0.07%  '
$ws.Range("E34").Value = '  -0.92%  '
$ws.Range("E35").Value = '  -3.29%  '
$ws.Range("E36").Value = '  -3.00%  '
$ws.Range("E37").Value = '  -2.54%  '
$ws.Range("E38").Value = '  -3.84%  '
$ws.Range("E39").Value = '  -1.38%  '
$ws.Range("E41").Value = '  +5.68%  '
$ws.Range("E42").Value = '  -1.74%  '
$ws.Range("E43").Value = '  -4.16%  '
$ws.Range("E44").Value = '  -2.68%  '
$ws.Range("E45").Value = '  +0.49%  '
$ws.Range("E46").Value = '  -0.59%  '
$ws.Range("E47").Value = '  +0.29%  '
$ws.Range("E48").Value = '  +3.61%  '
$ws.Range("E49").Value = '  -0.01%  '
$ws.Range("E50").Value = '  +12.26%  '
$ws.Range("E51").Value = '  -0.94%  '

# Row 46/47 swap: Coin name and Link columns
$ws.Range("B46").Value = 'Monero'
$ws.Range("B47").Value = 'VeChain'
$ws.Range("C46").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("C47").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
